$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 1: header row becomes ID, P1..P9 (was "Questões", 1..9) ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "P1"
$ws.Range("C1").Value = "P2"
$ws.Range("D1").Value = "P3"
$ws.Range("E1").Value = "P4"
$ws.Range("F1").Value = "P5"
$ws.Range("G1").Value = "P6"
$ws.Range("H1").Value = "P7"
$ws.Range("I1").Value = "P8"
$ws.Range("J1").Value = "P9"

# --- Row 2: fix the trailing period on G2 ---
$ws.Range("G2").Value = "Usei o R na monografia, mas não durante o mestrado"

# --- Formatting for header row: center/center alignment with explicit black font ---
# Build the combined style on an out-of-the-way helper cell first, then copy
# only the resulting format onto the header row so we don't leave behind a
# trail of intermediate/unused cell-style records.
$helper = $ws.Range("Z1")
$helper.HorizontalAlignment = -4108
$helper.VerticalAlignment = -4108
$helper.Font.Color = 0

$hdr = $ws.Range("A1:J1")
$helper.Copy()
$hdr.PasteSpecial(-4122)
$helper.Clear()

# --- Selection moves to G3 ---
$ws.Range("G3").Select()
